# Applies the edit described by the target diff to the active document.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First heading paragraph: add <w:ind w:firstLine="720"/>
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Format.FirstLineIndent = 36   # 36pt == 720 twips

# ---------------------------------------------------------------------
# 2) Locate the "Siehe Zusammenfassung : ..." list paragraph (last
#    non-empty paragraph before the trailing blank paragraph) and
#    replace it with a highlighted version (no bookmark - the bookmark
#    moves to the new final paragraph created below).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$pSiehe = $d.Paragraphs.Item($count - 1)

$xmlSiehe = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">Siehe Zusammenfassung : </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>„</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>Copy-Paste-Zusammenfassung Paper1</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>“</w:t></w:r></w:p></pkg:xmlData>
'@
$pSiehe.Range.InsertXML($xmlSiehe)

# ---------------------------------------------------------------------
# 3) Replace the trailing blank paragraph with the full block of new
#    content (new heading, new list items, and the final blank
#    paragraph that now carries the _GoBack bookmark).
# ---------------------------------------------------------------------
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)

$xmlBlock = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:p><w:pPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:pStyle w:val="berschrift1"/><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Protokoll Tagesleistungen 2</w:t></w:r><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>7</w:t></w:r><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>.11.2017</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">Siehe Zusammenfassung von Paper: A </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>framework</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>for</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> multimodal</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>…</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">T. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Upadhaya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> et. Al.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:ind w:left="1080"/><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>Ideen aus deren Erkenntnissen f</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>ü</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>r uns:</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>Statistisch genug Werte ben</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>ö</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">tigt um eine Aussage </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>ü</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>ber die Validit</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>ä</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">t </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>der Paper</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> machen zu k</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>ö</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>nnen (wahrscheinlich schwer, bei den wenigen [vergleichbaren] Daten)</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>Mehr als 2 MRT- Sequenzen vergleichen und statistisch interpretieren (Alleinstellungsmerkmal?)!</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">Den Autoren war die Standardisierung wichtig, also die Normierung der verschiedenen Aufnahmen. Kann man verschiedene Sequenzen auf die gleichen signifikanten Features dann </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>zur</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>ü</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>ck f</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>ü</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>hren</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>? Wahrscheinlich ein interessantes Teilgebiet zum Ausarbeiten f</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>ü</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="de-DE"/></w:rPr><w:t>r die SA</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:ind w:left="1080"/><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr></w:p>
</pkg:xmlData>
'@
$pLast.Range.InsertXML($xmlBlock)

# ---------------------------------------------------------------------
# 4) Make sure numId=3 (used by the new bullet-list paragraphs above)
#    resolves to a decimal "1)" numbered list, matching the new
#    abstractNum Word would mint when the user applies that list style.
# ---------------------------------------------------------------------
$pNum = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Statistisch genug*") {
        $pNum = $cand
        break
    }
}
if ($pNum -ne $null) {
    $pNum.Range.ListFormat.ApplyNumberDefault()
    $lt = $pNum.Range.ListFormat.ListTemplate
    $lvl = $lt.ListLevels.Item(1)
    $lvl.NumberFormat = "%1)"
}

Write-Host "done"
